$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "90.898.72"
$ws.Range("E2").Value = "  +1.48%  "
$ws.Range("D3").Value = "3.181.12"
$ws.Range("E3").Value = "  -1.61%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'215.00"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").Value = "'633.05"
$ws.Range("E6").Value = "  +1.85%  "
$ws.Range("D7").Value = "'0.406"
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").Value = "'0.730"
$ws.Range("E8").Value = "  +3.21%  "
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").Value = "3.181.40"
$ws.Range("E10").Value = "  -1.34%  "
$ws.Range("D11").Value = "'0.566"
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("E12").Value = "  +2.29%  "
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "90.631.45"
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("D15").Value = "'5.33"
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("D16").Value = "3.767.17"
$ws.Range("E16").Value = "  -1.15%  "
$ws.Range("D17").Value = "'32.49"
$ws.Range("E17").Value = "  -2.63%  "
$ws.Range("D18").Value = "3.187.33"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").Value = "'3.32"
$ws.Range("E19").Value = "  +2.57%  "
$ws.Range("E20").Value = "  +30.04%  "
$ws.Range("D21").Value = "'13.41"
$ws.Range("E21").Value = "  -3.23%  "
$ws.Range("D22").Value = "'433.98"
$ws.Range("E22").Value = "  +2.46%  "
$ws.Range("D23").Value = "'8.46"
$ws.Range("E23").Value = "  -2.72%  "
$ws.Range("D24").Value = "'4.99"
$ws.Range("E24").Value = "  -4.28%  "
$ws.Range("D25").Value = "'5.28"
$ws.Range("E25").Value = "  -2.73%  "
$ws.Range("D26").Value = "'11.68"
$ws.Range("E26").Value = "  -7.20%  "
$ws.Range("D27").Value = "'80.93"
$ws.Range("E27").Value = "  +7.94%  "
$ws.Range("D28").Value = "3.352.07"
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E30").Value = "  -8.80%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").Value = "'4.04"
$ws.Range("E32").Value = "  +25.06%  "
$ws.Range("D33").Value = "'8.36"
$ws.Range("E33").Value = "  -0.94%  "
$ws.Range("D34").Value = "'512.56"
$ws.Range("E34").Value = "  -9.34%  "
$ws.Range("D35").Value = "'6.95"
$ws.Range("E35").Value = "  -2.85%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Value = "'1.27"
$ws.Range("E37").Value = "  -4.85%  "
$ws.Range("D38").Value = "'22.31"
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("D39").Value = "'22.37"
$ws.Range("E39").Value = "  +2.20%  "
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("E41").Value = "  -5.22%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("E43").Value = "  -2.70%  "
$ws.Range("E44").Value = "  -2.92%  "
$ws.Range("D45").Value = "'147.52"
$ws.Range("E45").Value = "  -3.04%  "
$ws.Range("D46").Value = "'43.95"
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("D47").Value = "'169.09"
$ws.Range("E47").Value = "  -5.00%  "
$ws.Range("E48").Value = "  -1.59%  "
$ws.Range("D49").Value = "'0.739"
$ws.Range("E49").Value = "  +4.19%  "
$ws.Range("D50").Value = "'24.60"
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("E51").Value = "  -4.83%  "
